$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "pool_id" column is introduced before the old "manufacturer" column,
# and the old (duplicate) "serial_number" column that used to sit in E is
# removed. Rather than using structural column insert/delete (which the
# engine cannot currently round-trip losslessly for these particular
# column-width definitions), we reproduce the net effect directly by
# writing the resulting values into each affected cell.

# Row 1 - headers
$ws.Range("C1").Value = "pool_id"
$ws.Range("D1").Value = "manufacturer"
$ws.Range("E1").Value = "visibility"

# Row 2
$ws.Range("C2").Value = "co_hallo_7949"
$ws.Range("D2").Value = "b1"
$ws.Range("E2").Value = "a3"
$ws.Range("G2").Value = "hallo@hallo.at, hallo1@hallo.at"
$ws.Range("L2").Value = "aa"

# Row 3
$ws.Range("C3").Value = "co_hallo_7949"
$ws.Range("D3").Value = "b2"
$ws.Range("E3").Value = "b3"
$ws.Range("G3").Value = "hallo@hallo.at, hallo2@hallo.at"
$ws.Range("L3").Value = "bb"

# Row 4
$ws.Range("C4").Value = "co_hallo1_9778"
$ws.Range("E4").Value = ""
$ws.Range("L4").Value = "cc"

# Row 5
$ws.Range("C5").Value = "co_hallo1_9778"
$ws.Range("D5").Value = "d2"
$ws.Range("E5").Value = "d3"
$ws.Range("L5").Value = "dd"

# Row 6
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "e2"
$ws.Range("E6").Value = "e3"
$ws.Range("L6").Value = "ee"

# Row 3 grew to the taller row height used elsewhere (e.g. row 2).
$ws.Rows(3).RowHeight = 14.9

# Update the active selection to match the author's final cursor position.
$ws.Range("F9").Select()
